$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Turn on the table's Total Row (extends the table ref from G17 to G18,
# sets totalsRowCount="1" and marks totalsRowFunction="custom" on the
# columns that end up with a formula below).
$lo.ShowTotals = $true

# --- Totals row content ---------------------------------------------------
# A18: totals-row label
$ws.Range("A18").Value = "%"
$ws.Range("A18").HorizontalAlignment = -4108   # xlCenter
$ws.Range("A18").VerticalAlignment = -4108     # xlCenter
$ws.Range("A18").Borders.Item(7).LineStyle = 1 # xlEdgeLeft
$ws.Range("A18").Borders.Item(8).LineStyle = 1 # xlEdgeTop

# B18:F18 - percentage-of-affirmative-values formulas, same border style
foreach ($col in @("B", "C", "D", "E", "F")) {
    $addr = $col + "18"
    $srcRange = $col + "2:" + $col + "16"
    $ws.Range($addr).Formula = "=(SUM(" + $srcRange + ")/COUNT(" + $srcRange + "))*100"
    $cellRng = $ws.Range($addr)
    $cellRng.Borders.Item(7).LineStyle = 1   # xlEdgeLeft
    $cellRng.Borders.Item(10).LineStyle = 1  # xlEdgeRight
    $cellRng.Borders.Item(8).LineStyle = 1   # xlEdgeTop
}

# G18 - empty "comment" cell belonging to the totals row (no visible border)
$ws.Range("G18").Locked = $true

# Matches the selection left behind by the author after adding the row
$ws.Range("A18").Select()
